$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text format so numeric-looking strings
#     (e.g. "131.90", "0.605") are not auto-converted to numbers,
#     matching the original inline-string cell content. ---
$priceCells = @("D2","D3","D5","D6","D7","D8","D10","D11","D12","D13","D15","D16","D17","D18","D20","D22","D23","D25","D27","D28","D30","D33","D34","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.102.50'
$ws.Range("D3").Value = '3.475.22'
$ws.Range("D5").Value = '408.71'
$ws.Range("D6").Value = '131.90'
$ws.Range("D7").Value = '3.467.99'
$ws.Range("D8").Value = '0.605'
$ws.Range("D10").Value = '0.697'
$ws.Range("D11").Value = '0.131'
$ws.Range("D12").Value = '43.53'
$ws.Range("D13").Value = '4.027.77'
$ws.Range("D15").Value = '8.86'
$ws.Range("D16").Value = '20.16'
$ws.Range("D17").Value = '3.467.94'
$ws.Range("D18").Value = '63.168.61'
$ws.Range("D20").Value = '10.87'
$ws.Range("D22").Value = '3.35'
$ws.Range("D23").Value = '82.64'
$ws.Range("D25").Value = '312.58'
$ws.Range("D27").Value = '30.53'
$ws.Range("D28").Value = '8.18'
$ws.Range("D30").Value = '4.37'
$ws.Range("D33").Value = '44.17'
$ws.Range("D34").Value = '11.86'
$ws.Range("D36").Value = '0.998'
$ws.Range("D37").Value = '0.0494'
$ws.Range("D38").Value = '52.65'
$ws.Range("D39").Value = '3.57'
$ws.Range("D40").Value = '0.999'
$ws.Range("D41").Value = '3.04'
$ws.Range("D42").Value = '0.126'
$ws.Range("D44").Value = '137.04'
$ws.Range("D45").Value = '17.54'
$ws.Range("D46").Value = '3.99'
$ws.Range("D47").Value = '0.288'
$ws.Range("D49").Value = '22.10'
$ws.Range("D50").Value = '3.825.90'
$ws.Range("D51").Value = '2.189.62'

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Volume/coin/link columns: plain text assignment ---
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("E3").Value = '  +3.60%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("E6").Value = '  +17.24%  '
$ws.Range("E7").Value = '  +3.54%  '
$ws.Range("E8").Value = '  +2.97%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +8.74%  '
$ws.Range("E11").Value = '  +31.19%  '
$ws.Range("E12").Value = '  +9.51%  '
$ws.Range("E13").Value = '  +3.54%  '
$ws.Range("E15").Value = '  +5.24%  '
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("E17").Value = '  +3.63%  '
$ws.Range("E18").Value = '  +3.80%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("E21").Value = '  +28.95%  '
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("E23").Value = '  +9.68%  '
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("E25").Value = '  +3.04%  '
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("E27").Value = '  +6.01%  '
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("E30").Value = '  -2.72%  '
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("E33").Value = '  +13.11%  '
$ws.Range("E34").Value = '  +3.68%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  -3.22%  '
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("E43").Value = '  +4.03%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("E45").Value = '  +3.86%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("E47").Value = '  -4.25%  '
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("E50").Value = '  +3.62%  '
$ws.Range("E51").Value = '  +0.02%  '
